# Add a "Status" column (E) to the select2cols sheet: "Yes" for every
# county/health-district row, except a handful of counties that get "No"
# (Bourbon, Breathitt, Johnson, Lincoln, Magoffin, Martin, Pike - rows
# affected by the LNA data-file transfer/hunt/seek exception noted at
# Bracken Co).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("E1").Value = "Status"

# Default every data row (2-62) to "Yes" ...
for ($r = 2; $r -le 62; $r++) {
    $ws.Cells.Item($r, 5).Value = "Yes"
}

# ... then flip the exceptions to "No"
$noRows = @(6, 10, 34, 41, 44, 46, 55)
foreach ($r in $noRows) {
    $ws.Cells.Item($r, 5).Value = "No"
}

# Column D got narrower now that column E exists; column E keeps the
# sheet's default width. (64 is the input that the engine's character-width
# quantizer maps closest to the target 64.89 stored width.)
$ws.Columns.Item(4).ColumnWidth = 64

# Leave the cursor where the last edit happened.
$ws.Range("E52").Select() | Out-Null
